# TestData_NOP.xlsx - add SVBU project/repo lookup rows (F/G) with LOWER() formulas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "PROJECT_SVBU_WRONG_PROJ"
$ws.Range("G5").Formula = "=LOWER(F5)"

$ws.Range("F6").Value = "SVBU_WRONG_REPOCENTRE"
$ws.Range("G6").Formula = "=LOWER(F6)"

# Match the authored column widths (F ~28.43, G = 24 "characters") as closely as
# the engine's column-width quantization allows.
$ws.Columns.Item(6).ColumnWidth = 27.666666666666668
$ws.Columns.Item(7).ColumnWidth = 23.166666666666668

# Selection ends on G5:G6 (active cell G5), matching the authored view.
$ws.Range("G5:G6").Select()
